$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.06334924697876
$ws.Range("B1").Value = 5.190731525421143
$ws.Range("C1").Value = 3.279706239700317
$ws.Range("D1").Value = 2.292347431182861
$ws.Range("E1").Value = 2.090254545211792
